$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.632.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "'3.402.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'406.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").Value = "'130.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +16.61%  "
$ws.Range("D7").Value = "'0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.06%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.677"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.73%  "
$ws.Range("E10").Value = "  +11.27%  "
$ws.Range("D11").Value = "'42.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.34%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'3.950.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").Value = "'8.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "'19.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("D16").Value = "'3.393.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").Value = "'11.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.91%  "
$ws.Range("D18").Value = "'61.465.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("E19").Value = "  +4.88%  "
$ws.Range("D20").Value = "'0.0000135"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +18.38%  "
$ws.Range("D21").Value = "'3.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'82.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.52%  "
$ws.Range("D23").Value = "'13.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("D24").Value = "'308.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").Value = "'8.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.91%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'29.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.60%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'4.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.54%  "
$ws.Range("D29").Value = "'7.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("E31").Value = "  +5.15%  "
$ws.Range("D32").Value = "'11.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.78%  "
$ws.Range("E33").Value = "  +6.18%  "
$ws.Range("E34").Value = "  +10.44%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").Value = "'52.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "'3.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("D40").Value = "'3.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "'2.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.05%  "
$ws.Range("D42").Value = "'0.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.41%  "
$ws.Range("D43").Value = "'137.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("E44").Value = "  +8.35%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'17.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").Value = "'21.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("D49").Value = "'2.150.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'3.735.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("E51").Value = "  +0.32%  "
